$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Overall Demand") previously had empty cells in B2:D2; fill in the
# actual overall demand values for Iron & steel, Chemicals, and
# Non-metallic minerals respectively.
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 3612.120285859941
$ws.Range("D2").Value = 1499.915049436703
